# Helper: find a paragraph's 1-based index whose visible text (paragraph
# mark stripped) exactly equals $targetText.
function Get-ParaIndexByText {
    param($doc, [string]$targetText)
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $targetText) {
            return $i
        }
    }
    return -1
}

$d = $word.ActiveDocument

# --- "Where was I?" cell -------------------------------------------------
# Before:
#   De-server                     (ilvl 0)
#   Request bot                   (ilvl 1)
#   Contact page                  (ilvl 1)
#   Error messages in JS          (ilvl 1, carries the _GoBack bookmark)
#   Test with xampp stopped       (ilvl 1)
# After:
#   Images black and white        (ilvl 0, carries the _GoBack bookmark)
#
# The paragraph holding the bookmark ("Error messages in JS") is kept so the
# bookmark survives; its text becomes "Images black and white" and it is
# outdented to level 0. Every other paragraph in the block is removed.

$null = $d.Paragraphs.Item((Get-ParaIndexByText $d "Request bot")).Range.Delete()
$null = $d.Paragraphs.Item((Get-ParaIndexByText $d "Contact page")).Range.Delete()

$targetIdx = Get-ParaIndexByText $d "Error messages in JS"
$p = $d.Paragraphs.Item($targetIdx)
$r = $p.Range
$textOnly = $d.Range($r.Start, $r.End - 1)
$textOnly.Text = "Images black and white"
$d.Paragraphs.Item($targetIdx).Range.ListFormat.ListOutdent()

$null = $d.Paragraphs.Item((Get-ParaIndexByText $d "Test with xampp stopped")).Range.Delete()
$null = $d.Paragraphs.Item((Get-ParaIndexByText $d "De-server")).Range.Delete()

# --- "Next steps" cell -----------------------------------------------------
# The duplicate "Images black and white" bullet that used to sit right after
# "About page" is removed (the item now lives earlier in the document). Look
# for the occurrence immediately following "About page" specifically, since
# an earlier, different "Images black and white" paragraph now exists higher
# up in the document (the one created above) and must be left alone.
$aboutIdx = Get-ParaIndexByText $d "About page"
$dupIdx = $aboutIdx + 1
$dupText = $d.Paragraphs.Item($dupIdx).Range.Text.TrimEnd([char]13, [char]7)
if ($dupText -eq "Images black and white") {
    $d.Paragraphs.Item($dupIdx).Range.Delete()
} else {
    throw "Expected duplicate 'Images black and white' paragraph right after 'About page', found: [$dupText]"
}
